$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E (match bold style used by the other headers)
$ws.Range("E1").Value = "weight"
$ws.Range("E1").Font.Bold = $true

# New "weight" column values (rows 2-23)
$weights = @(
    79.599999999999994,
    79.8,
    79.3,
    79.2,
    79.3,
    79.400000000000006,
    79.599999999999994,
    79.8,
    79.8,
    79.400000000000006,
    79.3,
    79.400000000000006,
    79.2,
    79.099999999999994,
    79,
    78.900000000000006,
    78.8,
    78.7,
    78.8,
    78.599999999999994,
    78.5,
    78.599999999999994
)

for ($i = 0; $i -lt $weights.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $weights[$i]
}

# New row 23 - extra A value with no B/C
$ws.Range("A23").Value = 32

# Update selection to match the new active cell
$ws.Range("G11").Select()
